$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4574.91285437064
$ws.Range("D2").Value = 5016
$ws.Range("F2").Value = 3.27473414494826

$ws.Range("C3").Value = 4534.79981605076
$ws.Range("F3").Value = 111.603357549186

$ws.Range("C4").Value = 4218.13333459936
$ws.Range("F4").Value = 107.329841442977

$ws.Range("C5").Value = 1589.94976141188
$ws.Range("F5").Value = 20.1543227996469

$ws.Range("C6").Value = 1566.21463260172
$ws.Range("F6").Value = 18.7328125794779

$ws.Range("C7").Value = 3854.54701327964
$ws.Range("F7").Value = 100.650758892934

$ws.Range("C9").Value = 3786.9724574529
$ws.Range("F9").Value = 94.0596362612623

$ws.Range("C10").Value = 3669.73234755064
$ws.Range("F10").Value = 89.1746316820017

$ws.Range("C11").Value = 3609.82651413978
$ws.Range("F11").Value = 86.6785552898824

$ws.Range("C12").Value = 1273.83010908696
$ws.Range("F12").Value = 4.59402177271198

$ws.Range("C13").Value = 1257.73939713605
$ws.Range("F13").Value = 3.58491211011801

$ws.Range("C14").Value = 3539.86444065439
$ws.Range("F14").Value = 84.0811871727367

$ws.Range("C15").Value = 3536.05066810466
$ws.Range("F15").Value = 83.9222799831647
